# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: actualizar el resumen de conversión del día ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.14 = 7750.84 pesos`n✅ 7750.84 pesos = 2.14 = 961.11 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $text

# --- tasas: actualizar las tasas crudas usadas por las fórmulas ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 466.53
$ws2.Range("O10").Value = 3616
$ws2.Range("N12").Value = 3629
$ws2.Range("O12").Value = 450
